$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.072.56"
$ws.Range("E2").Value = "  +0.55%  "

$ws.Range("D3").Value = "2.312.29"
$ws.Range("E3").Value = "  +0.27%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.83%  "

$ws.Range("E7").Value = "  -0.34%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0914"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.96%  "

$ws.Range("E13").Value = "  -1.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.994"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.22%  "

$ws.Range("D16").Value = "2.663.62"

$ws.Range("D17").Value = "2.311.45"
$ws.Range("E17").Value = "  +0.36%  "

$ws.Range("D18").Value = "42.958.70"
$ws.Range("E18").Value = "  +0.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.35%  "

$ws.Range("E20").Value = "  -0.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.25%  "

$ws.Range("E25").Value = "  +0.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +15.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0872"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.43%  "

$ws.Range("E34").Value = "  +5.25%  "

$ws.Range("E35").Value = "  -1.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.112"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.37%  "

$ws.Range("E38").Value = "  +0.79%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.66"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "103.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.79%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.232"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "112.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.89%  "

$ws.Range("D48").Value = "1.661.35"
$ws.Range("E48").Value = "  -1.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.71%  "

$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.06%  "
